$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.182.25"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "'1.910.73"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D5").Value = "'314.84"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "'0.3923"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.09260"
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("D10").Value = "'1.139"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'41.93"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "'6.399"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "'20.88"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "'1.904.80"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'7.320"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'0.00001120"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "'92.37"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'0.06616"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'17.96"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'28.233.26"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "'2.330"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "'2.591"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").Value = "'2.126.97"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'158.17"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "'127.13"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "'1.103"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'5.640"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'3.612"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "'9.703"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").Value = "'0.06669"
$ws.Range("D37").Value = "'0.02428"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").Value = "'1.242"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "'0.2199"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "'1.286"
$ws.Range("E40").Value = "  +8.72%  "
$ws.Range("D41").Value = "'0.6476"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "'11.52"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'5.000"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'13.37"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "'0.6079"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").Value = "'1.291"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").Value = "'2.015"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "'123.55"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -0.74%  "
